$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.374.81'
$ws.Range('E2').Value = '  +7.37%  '

$ws.Range('D3').Value = '1.817.96'
$ws.Range('E3').Value = '  +7.90%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +1.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '344.29'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.77%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.89%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3831'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +5.12%  '

$ws.Range('B8').Value = 'OKB'
$ws.Range('C8').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '50.99'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +6.24%  '

$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3505'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +7.05%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.226'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +6.05%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07738'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +6.00%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.002'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.11%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.02'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +10.91%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.613'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +7.88%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.234'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +7.46%  '

$ws.Range('D16').Value = '1.818.91'
$ws.Range('E16').Value = '  +8.49%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001120'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +5.75%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06757'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +3.36%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '87.05'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +9.05%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.89%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.86'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +11.68%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.507'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +9.00%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.11'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.78%  '

$ws.Range('D24').Value = '27.414.92'
$ws.Range('E24').Value = '  +7.63%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.467'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.91%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.685'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +9.57%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.00'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +16.33%  '

$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.501'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +17.21%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '153.04'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.09%  '

$ws.Range('D30').Value = '2.022.03'
$ws.Range('E30').Value = '  +8.51%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '137.04'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +7.62%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.186'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +2.80%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.313'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +5.95%  '

$ws.Range('B34').Value = 'Stellar'
$ws.Range('C34').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08812'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +4.05%  '

$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '13.62'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +8.50%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.727'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +3.54%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.642'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +7.56%  '

$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '9.099'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +8.18%  '

$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06569'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +7.20%  '

$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02424'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +7.75%  '

$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2255'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +7.46%  '

$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6835'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +13.52%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.253'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.60%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.86'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +8.14%  '

$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.000'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.93%  '

$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6397'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +10.79%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.970'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +3.99%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.179'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +9.92%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '132.66'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +6.41%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07369'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.39%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '80.74'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +7.38%  '
